$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "986÷3=328, 2"; New = "121÷2=60, 1" },
    @{ Old = "120÷7=17, 1"; New = "867÷9=96, 3" },
    @{ Old = "848÷5=169, 3"; New = "326÷8=40, 6" },
    @{ Old = "512÷8=64, 0"; New = "199÷2=99, 1" },
    @{ Old = "563÷4=140, 3"; New = "115÷7=16, 3" },
    @{ Old = "185÷8=23, 1"; New = "684÷9=76, 0" },
    @{ Old = "438÷8=54, 6"; New = "890÷3=296, 2" },
    @{ Old = "391÷2=195, 1"; New = "602÷7=86, 0" },
    @{ Old = "419÷5=83, 4"; New = "963÷5=192, 3" },
    @{ Old = "294÷6=49, 0"; New = "321÷8=40, 1" },
    @{ Old = "620÷4=155, 0"; New = "228÷5=45, 3" },
    @{ Old = "859÷3=286, 1"; New = "352÷8=44, 0" },
    @{ Old = "523÷3=174, 1"; New = "603÷9=67, 0" },
    @{ Old = "832÷7=118, 6"; New = "350÷9=38, 8" },
    @{ Old = "833÷7=119, 0"; New = "827÷7=118, 1" },
    @{ Old = "914÷9=101, 5"; New = "545÷6=90, 5" },
    @{ Old = "834÷8=104, 2"; New = "551÷4=137, 3" },
    @{ Old = "220÷5=44, 0"; New = "610÷6=101, 4" },
    @{ Old = "579÷2=289, 1"; New = "140÷5=28, 0" },
    @{ Old = "642÷2=321, 0"; New = "579÷5=115, 4" },
    @{ Old = "328÷9=36, 4"; New = "592÷9=65, 7" },
    @{ Old = "475÷3=158, 1"; New = "163÷2=81, 1" },
    @{ Old = "231÷6=38, 3"; New = "911÷8=113, 7" },
    @{ Old = "214÷9=23, 7"; New = "594÷8=74, 2" },
    @{ Old = "430÷4=107, 2"; New = "281÷9=31, 2" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.New, 2)
}
